$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Landing Page").Name = "Quotes"
$wb.Worksheets.Item("reading").Name = "readingQuotes"
$wb.Worksheets.Item("inspirational").Name = "inspirationalQuotes"
$wb.Worksheets.Item("friendship").Name = "friendshipQuotes"
